# Replace the sample patient's record with a new patient's record
# (commit: "para agregar codigo de barras").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient name / expediente (row 6) ---
$ws.Range("A6").Value = "MORAN "
$ws.Range("C6").Value = "VARGAS"
$ws.Range("E6").Value = "IRIS"
$ws.Range("G6").Value = "AMABELY"
$ws.Range("I6").Value = "/201773459"

# --- Dirección actual (row 8: calle, municipio, departamento, telefono) ---
$ws.Range("A8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("F8").Value = "SANTA ELENA"
$ws.Range("H8").Value = "PETEN"
$ws.Range("J8").Value = "4618-2848"

# --- Fecha de nacimiento / edad / lugar de nacimiento (row 12) ---
# Force text format so the date-like / numeric-like strings are not
# auto-converted to a date serial / number by Excel.
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "1991-05-12"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "26"
$ws.Range("H12").Value = "PETEN"

# --- Estado civil (row 14) ---
$ws.Range("A14").Value = "Soltero"

# --- No. de Cédula (row 14) ---
$ws.Range("H14").Value = "SIN DOC."

# --- Nombre del Cónyugue (row 16) ---
$ws.Range("A16").Value = "MELVIN LÓPEZ"

# --- Nombre del Padre / Nombre de la Madre (row 18) ---
$ws.Range("A18").Value = ""
$ws.Range("F18").Value = "NORA VARGAS"

# --- En caso de emergencia: nombre / parentesco / telefono (row 20) ---
$ws.Range("A20").Value = "MELVIN LÓPEZ"
$ws.Range("F20").Value = "ESPOSO"
$ws.Range("J20").Value = "3134-2545"

# --- Fecha de ingreso / hora / servicio (row 24) ---
$ws.Range("A24").Value = "20/11/2017"
$ws.Range("C24").Value = "14:29:3"
$ws.Range("D24").Value = "UNIDAD 18"
